# Updates the cryptos price/volume table to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.767.64'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.628.76'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  -0.85%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.85'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0630'
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.57'
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0787'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.854.98'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.581.98'
$ws.Range('E14').Value = '  -3.44%  '
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0756'
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.66'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.784.98'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.42'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.71'
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.22'
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('E27').Value = '  +1.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.82'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.49'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.21'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.58'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.901'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.140.43'
$ws.Range('E37').Value = '  +2.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.543'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.995'
$ws.Range('E41').Value = '  -0.81%  '
$ws.Range('E42').Value = '  -1.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.60'
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.30'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.797'
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.766.38'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.32'
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.46'
$ws.Range('E49').Value = '  +6.27%  '
$ws.Range('E50').Value = '  +2.02%  '
$ws.Range('E51').Value = '  -0.54%  '
